# "Improve VGA control code"
# Add a new time-record entry for 1.4.2020 (10:00-10:30, VGA Control /
# Improve RTL) to the bottom of the existing task list, and leave the
# sheet scrolled to the top with the new row's last cell selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the look (alignment / number formats) of the row above - it belongs
# to the same task block - onto the new row before filling in its data.
$ws.Range("A17:F17").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)  # xlPasteFormats

# A: date, stored as plain text like the rest of column A (e.g. "31.3.2020"),
# not as a real date serial. Go through a text formula then paste its
# value back so Excel doesn't reinterpret the "1.4.2020" string as a date.
$ws.Range("A19").Formula = "=""1.4.2020"""
$ws.Range("A19").Copy()
$ws.Range("A19").PasteSpecial(-4163)  # xlPasteValues

# B/C: From 10:00 to 10:30
$ws.Range("B19").Value = 0.41666666666666669
$ws.Range("C19").Value = 0.4375

# D: worked time = To - From
$ws.Range("D19").Formula = "=C19-B19"

# E/F: task / notes
$ws.Range("E19").Value = "VGA Control"
$ws.Range("F19").Value = "Improve RTL"

# Scroll back to the top of the sheet and select the newly-entered note cell.
$ws.Range("A1").Select()
$ws.Range("F19").Select()
